$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("M24").Value = 9735.9
$wsVentasPorGrupo.Range("M25").Value = "4 de 23"

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F24").Value = 9735.9
$wsVentaMensual.Range("F25").Value = 30206.69

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 29393.57
$wsCumplimiento.Range("E12").Value = 13706.5154117774
$wsCumplimiento.Range("F12").Value = 0.6819840313348429

$wsCumplimiento.Range("D15").Value = 30206.69
$wsCumplimiento.Range("E15").Value = 27996.77623249458
$wsCumplimiento.Range("F15").Value = 0.5189843828087307
